# Insert 2 new rows right after the header row, shifting the existing 20 data rows down
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:3").Insert()
$ws.Range("A2:H3").ClearFormats()

# Fix up the timestamp (column A) sequence for every data row now that rows 4-23 shifted down by two,
# so the whole A column remains a clean 0,100,200,... ladder through row 23.
For ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
    $ws.Cells.Item($r, 2).Value = "struggle"
}

# Populate the two brand-new rows inserted at the top of the data (rows 2-3)
$ws.Range("C2").Value = 0.001320064067840854
$ws.Range("D2").Value = 0.1883212360553446
$ws.Range("E2").Value = -0.001869207888376141
$ws.Range("F2").Value = 0.5199990272521973
$ws.Range("G2").Value = -4.081954479217529
$ws.Range("H2").Value = 1.088714599609375

$ws.Range("C3").Value = -0.1889566183090211
$ws.Range("D3").Value = 0.01437168661504978
$ws.Range("E3").Value = 0.113000919460319
$ws.Range("F3").Value = 0.446847915649414
$ws.Range("G3").Value = -3.855323314666748
$ws.Range("H3").Value = 1.094059705734253

# Append 8 brand-new rows at the bottom of the data (rows 24-31)
$ws.Range("A24").Value = 2200
$ws.Range("B24").Value = "struggle"
$ws.Range("C24").Value = 4.677844420075353
$ws.Range("D24").Value = -3.651133604347696
$ws.Range("E24").Value = -7.842656075954431
$ws.Range("F24").Value = -0.7247915863990784
$ws.Range("G24").Value = -2.964529037475586
$ws.Range("H24").Value = -2.036930084228516

$ws.Range("A25").Value = 2300
$ws.Range("B25").Value = "struggle"
$ws.Range("C25").Value = -2.627654522657398
$ws.Range("D25").Value = -2.928949266672134
$ws.Range("E25").Value = 4.230176210403448
$ws.Range("F25").Value = -0.9755517840385436
$ws.Range("G25").Value = -3.013092756271362
$ws.Range("H25").Value = -0.1954768747091293

$ws.Range("A26").Value = 2400
$ws.Range("B26").Value = "struggle"
$ws.Range("C26").Value = -4.852406792342663
$ws.Range("D26").Value = 0.3913787733763447
$ws.Range("E26").Value = 0.2968738228082666
$ws.Range("F26").Value = -0.0148134818300604
$ws.Range("G26").Value = -4.330729007720947
$ws.Range("H26").Value = 0.6641632318496704

$ws.Range("A27").Value = 2500
$ws.Range("B27").Value = "struggle"
$ws.Range("C27").Value = -1.301035702228551
$ws.Range("D27").Value = 3.64691380783915
$ws.Range("E27").Value = -6.109266191720954
$ws.Range("F27").Value = -0.2535090744495392
$ws.Range("G27").Value = -4.50192403793335
$ws.Range("H27").Value = 0.8677340745925903

$ws.Range("A28").Value = 2600
$ws.Range("B28").Value = "struggle"
$ws.Range("C28").Value = 2.465943455696097
$ws.Range("D28").Value = -2.991184197366218
$ws.Range("E28").Value = -3.608212560415278
$ws.Range("F28").Value = -0.0308486949652433
$ws.Range("G28").Value = -3.680310487747192
$ws.Range("H28").Value = 1.009607553482056

$ws.Range("A29").Value = 2700
$ws.Range("B29").Value = "struggle"
$ws.Range("C29").Value = -1.307898223400096
$ws.Range("D29").Value = -2.068972408771528
$ws.Range("E29").Value = -0.7334359884262174
$ws.Range("F29").Value = -0.2924517393112182
$ws.Range("G29").Value = 0.6568328738212585
$ws.Range("H29").Value = 0.4216497242450714

$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "struggle"
$ws.Range("C30").Value = -1.702915767207749
$ws.Range("D30").Value = -0.5735956337302961
$ws.Range("E30").Value = -0.9715757742524092
$ws.Range("F30").Value = 0.683863639831543
$ws.Range("G30").Value = 4.383111000061035
$ws.Range("H30").Value = -1.505782842636108

$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "struggle"
$ws.Range("C31").Value = -3.414293382316824
$ws.Range("D31").Value = 0.2869436666369428
$ws.Range("E31").Value = -0.1008520126342796
$ws.Range("F31").Value = 0.6068946123123169
$ws.Range("G31").Value = 4.862334728240967
$ws.Range("H31").Value = -0.4990769028663635

